$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric-looking "Price" values are stored as text in this workbook.
# Writing them with a leading apostrophe via .Formula forces Excel to keep
# them as text instead of silently converting to a Number, and resetting
# the cell Style back to "Normal" afterwards clears the quote-prefix style
# that operation would otherwise leave behind.
$ws.Range("D2").Formula = "'262.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Formula = "'22.97"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Formula = "'6.190"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Formula = "'0.06096"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Formula = "'3.463"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Formula = "'1.368"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Formula = "'0.7994"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Formula = "'0.08119"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Formula = "'0.03493"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Formula = "'0.03046"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Formula = "'0.09319"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Formula = "'3.854"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Formula = "'0.001702"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Formula = "'0.04794"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Formula = "'0.0006140"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Formula = "'0.006206"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Formula = "'0.001094"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Formula = "'0.003807"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Formula = "'3.703"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Formula = "'2.215"
$ws.Range("D24").Style = "Normal"
$ws.Range("D26").Formula = "'0.1253"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Formula = "'0.0003202"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Formula = "'0.04603"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Formula = "'0.007060"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Formula = "'0.003900"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Formula = "'0.1118"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Formula = "'0.01069"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Formula = "'0.002970"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Formula = "'0.00005938"
$ws.Range("D46").Style = "Normal"
$ws.Range("D48").Formula = "'0.7000"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Formula = "'0.07326"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Formula = "'0.00002100"
$ws.Range("D50").Style = "Normal"

# Plain text fields (coin name / link / label) - swap BKEXToken <-> KickToken rows
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E43").Value = "42BKEXTokenBKK"
